{"js": "// Update the worksheet date and the 25 division problems (5 rows of 5\n// problems, each followed by 3 blank rows for student work) to the new\n// day's values. Cells are addressed positionally (row, col) rather than\n// by searching for their current text, because several of the problems\n// share the same \"before\" text (e.g. \"45\u00f73=\" and \"66\u00f77=\" each occur\n// twice) and a text-search replace-all would not be able to tell them\n// apart.\n\n// 1) Header date paragraph: \"2023-10-16 Monday\" -> \"2023-10-17 Tuesday\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items.find(p => p.text.trim() === \"2023-10-16 Monday\");\nif (dateParagraph) {\n  dateParagraph.getRange().insertText(\"2023-10-17 Tuesday\", \"Replace\");\n}\n\n// 2) The division problems table: 5 \"problem\" rows (0, 4, 8, 12, 16),\n// each with 5 columns, interleaved with blank work rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"56\u00f75=\" },\n  { row: 0, col: 1, text: \"56\u00f74=\" },\n  { row: 0, col: 2, text: \"92\u00f74=\" },\n  { row: 0, col: 3, text: \"30\u00f73=\" },\n  { row: 0, col: 4, text: \"10\u00f73=\" },\n\n  { row: 4, col: 0, text: \"62\u00f74=\" },\n  { row: 4, col: 1, text: \"14\u00f79=\" },\n  { row: 4, col: 2, text: \"52\u00f72=\" },\n  { row: 4, col: 3, text: \"22\u00f75=\" },\n  { row: 4, col: 4, text: \"97\u00f72=\" },\n\n  { row: 8, col: 0, text: \"90\u00f73=\" },\n  { row: 8, col: 1, text: \"25\u00f75=\" },\n  { row: 8, col: 2, text: \"65\u00f72=\" },\n  { row: 8, col: 3, text: \"51\u00f79=\" },\n  { row: 8, col: 4, text: \"51\u00f76=\" },\n\n  { row: 12, col: 0, text: \"38\u00f79=\" },\n  { row: 12, col: 1, text: \"54\u00f72=\" },\n  { row: 12, col: 2, text: \"73\u00f74=\" },\n  { row: 12, col: 3, text: \"90\u00f75=\" },\n  { row: 12, col: 4, text: \"69\u00f73=\" },\n\n  { row: 16, col: 0, text: \"95\u00f73=\" },\n  { row: 16, col: 1, text: \"37\u00f79=\" },\n  { row: 16, col: 2, text: \"35\u00f77=\" },\n  { row: 16, col: 3, text: \"80\u00f77=\" },\n  { row: 16, col: 4, text: \"60\u00f79=\" },\n];\n\nfor (const { row, col, text } of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division problems (5 rows of 5\n# problems, each followed by 3 blank rows for student work) to the new\n# day's values. Table cells are addressed positionally (row, col) rather\n# than by searching for their current text, because several of the\n# problems share the same \"before\" text (e.g. \"45\u00f73=\" and \"66\u00f77=\" each\n# occur twice) and a text search/replace-all would not be able to tell\n# them apart.\n\n$d = $word.ActiveDocument\n\n# 1) Header date paragraph: \"2023-10-16 Monday\" -> \"2023-10-17 Tuesday\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\"2023-10-16 Monday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-10-17 Tuesday\", 2) | Out-Null\n\n# 2) The division problems table: 5 \"problem\" rows (1, 5, 9, 13, 17 in\n# Word's 1-based row numbering), each with 5 columns, interleaved with\n# blank work rows.\n$table = $d.Tables.Item(1)\n\n$updates = @(\n  @{ Row = 1;  Col = 1; Text = \"56\u00f75=\" },\n  @{ Row = 1;  Col = 2; Text = \"56\u00f74=\" },\n  @{ Row = 1;  Col = 3; Text = \"92\u00f74=\" },\n  @{ Row = 1;  Col = 4; Text = \"30\u00f73=\" },\n  @{ Row = 1;  Col = 5; Text = \"10\u00f73=\" },\n\n  @{ Row = 5;  Col = 1; Text = \"62\u00f74=\" },\n  @{ Row = 5;  Col = 2; Text = \"14\u00f79=\" },\n  @{ Row = 5;  Col = 3; Text = \"52\u00f72=\" },\n  @{ Row = 5;  Col = 4; Text = \"22\u00f75=\" },\n  @{ Row = 5;  Col = 5; Text = \"97\u00f72=\" },\n\n  @{ Row = 9;  Col = 1; Text = \"90\u00f73=\" },\n  @{ Row = 9;  Col = 2; Text = \"25\u00f75=\" },\n  @{ Row = 9;  Col = 3; Text = \"65\u00f72=\" },\n  @{ Row = 9;  Col = 4; Text = \"51\u00f79=\" },\n  @{ Row = 9;  Col = 5; Text = \"51\u00f76=\" },\n\n  @{ Row = 13; Col = 1; Text = \"38\u00f79=\" },\n  @{ Row = 13; Col = 2; Text = \"54\u00f72=\" },\n  @{ Row = 13; Col = 3; Text = \"73\u00f74=\" },\n  @{ Row = 13; Col = 4; Text = \"90\u00f75=\" },\n  @{ Row = 13; Col = 5; Text = \"69\u00f73=\" },\n\n  @{ Row = 17; Col = 1; Text = \"95\u00f73=\" },\n  @{ Row = 17; Col = 2; Text = \"37\u00f79=\" },\n  @{ Row = 17; Col = 3; Text = \"35\u00f77=\" },\n  @{ Row = 17; Col = 4; Text = \"80\u00f77=\" },\n  @{ Row = 17; Col = 5; Text = \"60\u00f79=\" }\n)\n\nforeach ($u in $updates) {\n  $cell = $table.Cell($u.Row, $u.Col)\n  $cell.Range.Text = $u.Text\n}\n"}
